$d = $word.ActiveDocument
$wdNS = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Get-ParaByText($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text.StartsWith($needle)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Change 1: "Bill Fisher, Ron Richardson, Christina  Eusanio" paragraph
#   - merge the " " + "Eusanio" runs (dropping the spell-check proofErr tags)
#   - move the "_GoBack" bookmark here (it used to sit mid-paragraph 7)
# ---------------------------------------------------------------------------
$pTeam = Get-ParaByText("Bill Fisher")
$xmlTeam = '<w:p xmlns:w="' + $wdNS + '">' +
    '<w:r><w:t>Bill Fisher, Ron Richardson, Christina</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Eusanio</w:t></w:r>' +
    '<w:bookmarkStart w:id="100" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="100"/>' +
    '</w:p>'
$pTeam.Range.InsertXML($xmlTeam)

# ---------------------------------------------------------------------------
# Change 2: rewrite + expand the "Reddit is the fifth most popular..." para
#   into 3 paragraphs (the old mid-paragraph bookmark is gone, replaced above)
# ---------------------------------------------------------------------------
$pReddit = Get-ParaByText("Reddit is the fifth most popular")
$para1 = '<w:p xmlns:w="' + $wdNS + '">' +
    '<w:r><w:t xml:space="preserve">Predicting </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">page popularity, no matter what the platform, has many practical uses.  Namely, advertisers and marketers can leverage this information to target certain posts, pages, etc.  </w:t></w:r>' +
    '<w:r><w:t>Reddit is the fifth most popular site on the internet according to</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> the data captured by</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Amazon' + [char]0x2019 + 's Alexa.  Reddit is a social message board where users can create profiles, create topics, and interact with others in a variety of ways.  </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">When users create posts, they receive upvotes and downvotes.  The aggregate of those votes constitutes as a post' + [char]0x2019 + 's score.  </w:t></w:r>' +
    '<w:r><w:t>By extracting</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> various features from the Reddit API</w:t></w:r>' +
    '<w:r><w:t>, our goal is</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> to create a machine learning model that will </w:t></w:r>' +
    '<w:r><w:t>classify posts as popular or unpopular, as defined below</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '</w:p>'

$para2 = '<w:p xmlns:w="' + $wdNS + '">' +
    '<w:r><w:t xml:space="preserve">Our plan is to use the framework from </w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t>Popularity Prediction of Reddit Texts</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">, a thesis paper by </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">Tracy </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Rohlin</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">.  The author narrowed down her </w:t></w:r>' +
    '<w:r><w:t>dataset to 1</w:t></w:r>' +
    '<w:r><w:t>2</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">,000 posts from six subreddits.  We are going to focus on the following subreddits: </w:t></w:r>' +
    '<w:r><w:t>r/aww, r/fitness, r/conspiracy, r/knitting, r/</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>askmen</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>,</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> and</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> r/</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>askwomen</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">.  Instead of splitting the entire dataset into train, validation, and test sets, the author </w:t></w:r>' +
    '<w:r><w:t>deliberately</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> broke up each subreddit into train, validation, and test sets</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> that were stratified to include the same ratio of popular and unpopular posts within each</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">.  This way each subreddits' + [char]0x2019 + ' posts would be evaluated by the model equally.  The author also purposely uses posts that are roughly six months old because at that point voting has ended.  </w:t></w:r>' +
    '<w:r><w:t>Finally, the author labeled the data as popular or unpopular.  Posts were deemed popular if they fell above the 75</w:t></w:r>' +
    '<w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> percentile of voting scores</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> from the collected dataset</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '</w:p>'

$para3 = '<w:p xmlns:w="' + $wdNS + '">' +
    '<w:r><w:t>The author preprocessed the data by removing stop words, low frequency words (those that only show up in &lt;= two posts), and, punctuation.  The author also chose not to use stemming or lemmatization.  The posts were th</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">en fed to various preprocessing algorithms before being fed into predictor algorithms.  Finally, the author uses F1 scores and accuracy to evaluate the algorithms.  </w:t></w:r>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $wdNS + '"/>'

$pReddit.Range.InsertXML($para1 + $para2 + $para3)

# ---------------------------------------------------------------------------
# Change 3: "Sub Reddit" / "Sub Reddit data" -> split into Sub / - / Reddit[...]
# ---------------------------------------------------------------------------
$pSub1 = Get-ParaByText("Sub Reddit")
while ($pSub1 -ne $null -and $pSub1.Range.Text.TrimEnd() -ne "Sub Reddit") {
    $pSub1 = $null
}
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq "Sub Reddit") {
        $xml = '<w:p xmlns:w="' + $wdNS + '">' +
            '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
            '<w:r><w:t>Sub</w:t></w:r>' +
            '<w:r><w:t>-</w:t></w:r>' +
            '<w:r><w:t>Reddit</w:t></w:r>' +
            '</w:p>'
        $p.Range.InsertXML($xml)
        break
    }
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq "Sub Reddit data") {
        $xml = '<w:p xmlns:w="' + $wdNS + '">' +
            '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
            '<w:r><w:t>Sub</w:t></w:r>' +
            '<w:r><w:t>-</w:t></w:r>' +
            '<w:r><w:t>Reddit data</w:t></w:r>' +
            '</w:p>'
        $p.Range.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# Change 4: add <w:lastRenderedPageBreak/> before "Create date/time"
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq "Create date/time") {
        $xml = '<w:p xmlns:w="' + $wdNS + '">' +
            '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
            '<w:r><w:lastRenderedPageBreak/><w:t>Create date/time</w:t></w:r>' +
            '</w:p>'
        $p.Range.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# Change 5: add a new bullet "PRAW (Python Reddit API Wrapper)" right after
#   the "Scikit Learn" bullet
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq "Scikit Learn") {
        $r = $p.Range
        $r.Collapse(0) | Out-Null
        $xml = '<w:p xmlns:w="' + $wdNS + '">' +
            '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
            '<w:r><w:t>PRAW (Python Reddit API Wrapper)</w:t></w:r>' +
            '</w:p>'
        $r.InsertXML($xml)
        break
    }
}

Write-Host "Done applying edits"
